$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 313
$ws.Range("C2").Value = 313
$ws.Range("D2").Value = 300
$ws.Range("E2").Value = 300
$ws.Range("B3").Value = 397
$ws.Range("C3").Value = 467
$ws.Range("D3").Value = 300
$ws.Range("E3").Value = 300
$ws.Range("F3").Value = 75900
$ws.Range("G3").Value = 84600
$ws.Range("B4").Value = 279
$ws.Range("C4").Value = 494
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 79400
$ws.Range("G4").Value = 93400
$ws.Range("B5").Value = 406
$ws.Range("C5").Value = 485
$ws.Range("D5").Value = 400
$ws.Range("E5").Value = 400
$ws.Range("F5").Value = 111600
$ws.Range("G5").Value = 197600
$ws.Range("B6").Value = 579
$ws.Range("C6").Value = 668
$ws.Range("D6").Value = 800
$ws.Range("E6").Value = 800
$ws.Range("F6").Value = 148800
$ws.Range("G6").Value = 187200
$ws.Range("B7").Value = 423
$ws.Range("C7").Value = 243
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = 300
$ws.Range("F7").Value = 115800
$ws.Range("G7").Value = 200400
$ws.Range("B8").Value = 404
$ws.Range("C8").Value = 320
$ws.Range("D8").Value = 400
$ws.Range("E8").Value = 600
$ws.Range("F8").Value = 169200
$ws.Range("G8").Value = 145800
$ws.Range("B9").Value = 436
$ws.Range("C9").Value = 647
$ws.Range("D9").Value = 800
$ws.Range("E9").Value = 1200
$ws.Range("F9").Value = 227200
$ws.Range("G9").Value = 126000
$ws.Range("B10").Value = 267
$ws.Range("C10").Value = 444
$ws.Range("D10").Value = 300
$ws.Range("E10").Value = 300
$ws.Range("F10").Value = 130800
$ws.Range("G10").Value = 194100
$ws.Range("B11").Value = 461
$ws.Range("C11").Value = 287
$ws.Range("D11").Value = 600
$ws.Range("E11").Value = 600
$ws.Range("F11").Value = 145200
$ws.Range("G11").Value = 240000
$ws.Range("B12").Value = 498
$ws.Range("C12").Value = 374
$ws.Range("D12").Value = 400
$ws.Range("E12").Value = 500
$ws.Range("F12").Value = 184400
$ws.Range("G12").Value = 143500
$ws.Range("B13").Value = 658
$ws.Range("C13").Value = 632
$ws.Range("D13").Value = 800
$ws.Range("E13").Value = 1000
$ws.Range("F13").Value = 225600
$ws.Range("G13").Value = 149000
$ws.Range("B14").Value = 331
$ws.Range("C14").Value = 445
$ws.Range("D14").Value = 200
$ws.Range("E14").Value = 400
$ws.Range("F14").Value = 130200
$ws.Range("G14").Value = 252800
$ws.Range("B15").Value = 595
$ws.Range("C15").Value = 682
$ws.Range("D15").Value = 700
$ws.Range("E15").Value = 800
$ws.Range("F15").Value = 156800
$ws.Range("G15").Value = 196000
$ws.Range("B16").Value = 479
$ws.Range("C16").Value = 445
$ws.Range("D16").Value = 300
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 167100
$ws.Range("G16").Value = 191100
$ws.Range("B17").Value = 373
$ws.Range("C17").Value = 445
$ws.Range("D17").Value = 500
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 192000
$ws.Range("G17").Value = 172400
$ws.Range("B18").Value = 327
$ws.Range("C18").Value = 480
$ws.Range("D18").Value = 500
$ws.Range("E18").Value = 600
$ws.Range("F18").Value = 148500
$ws.Range("G18").Value = 174600
$ws.Range("B19").Value = 379
$ws.Range("C19").Value = 549
$ws.Range("D19").Value = 600
$ws.Range("E19").Value = 400
$ws.Range("F19").Value = 114600
$ws.Range("G19").Value = 172400
$ws.Range("B20").Value = 272
$ws.Range("C20").Value = 301
$ws.Range("D20").Value = 500
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 157500
$ws.Range("G20").Value = 206400
$ws.Range("B21").Value = 323
$ws.Range("C21").Value = 261
$ws.Range("D21").Value = 800
$ws.Range("E21").Value = 700
$ws.Range("F21").Value = 138400
$ws.Range("G21").Value = 174300
$ws.Range("B22").Value = 409
$ws.Range("C22").Value = 261
$ws.Range("D22").Value = 600
$ws.Range("E22").Value = 700
$ws.Range("F22").Value = 172200
$ws.Range("G22").Value = 151200
$ws.Range("B23").Value = 398
$ws.Range("C23").Value = 411
$ws.Range("D23").Value = 500
$ws.Range("E23").Value = 700
$ws.Range("F23").Value = 202000
$ws.Range("G23").Value = 178500
$ws.Range("B24").Value = 286
$ws.Range("C24").Value = 530
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 700
$ws.Range("F24").Value = 173500
$ws.Range("G24").Value = 163100
$ws.Range("B25").Value = 300
$ws.Range("C25").Value = 469
$ws.Range("D25").Value = 600
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 151800
$ws.Range("G25").Value = 195200
$ws.Range("B26").Value = 459
$ws.Range("C26").Value = 431
$ws.Range("D26").Value = 800
$ws.Range("E26").Value = 500
$ws.Range("F26").Value = 141600
$ws.Range("G26").Value = 185000
$ws.Range("B27").Value = 282
$ws.Range("C27").Value = 433
$ws.Range("D27").Value = 400
$ws.Range("E27").Value = 500
$ws.Range("F27").Value = 164800
$ws.Range("G27").Value = 157000
$ws.Range("B28").Value = 378
$ws.Range("C28").Value = 472
$ws.Range("D28").Value = 700
$ws.Range("E28").Value = 500
$ws.Range("F28").Value = 161700
$ws.Range("G28").Value = 191000
$ws.Range("B29").Value = 405
$ws.Range("C29").Value = 303
$ws.Range("D29").Value = 500
$ws.Range("E29").Value = 500
$ws.Range("F29").Value = 175500
$ws.Range("G29").Value = 230500
$ws.Range("B30").Value = 633
$ws.Range("C30").Value = 507
$ws.Range("D30").Value = 600
$ws.Range("E30").Value = 800
$ws.Range("F30").Value = 156000
$ws.Range("G30").Value = 148800
$ws.Range("B31").Value = 301
$ws.Range("C31").Value = 591
$ws.Range("D31").Value = 300
$ws.Range("E31").Value = 500
$ws.Range("F31").Value = 163200
$ws.Range("G31").Value = 176500
$ws.Range("B32").Value = 493
$ws.Range("C32").Value = 459
$ws.Range("D32").Value = 600
$ws.Range("E32").Value = 400
$ws.Range("F32").Value = 135000
$ws.Range("G32").Value = 212800
$ws.Range("B33").Value = 618
$ws.Range("C33").Value = 399
$ws.Range("D33").Value = 400
$ws.Range("E33").Value = 500
$ws.Range("F33").Value = 146800
$ws.Range("G33").Value = 169500
$ws.Range("B34").Value = 237
$ws.Range("C34").Value = 396
$ws.Range("D34").Value = 300
$ws.Range("E34").Value = 500
$ws.Range("F34").Value = 185400
$ws.Range("G34").Value = 177500
$ws.Range("B35").Value = 284
$ws.Range("C35").Value = 240
$ws.Range("D35").Value = 600
$ws.Range("E35").Value = 600
$ws.Range("F35").Value = 141000
$ws.Range("G35").Value = 237600
$ws.Range("B36").Value = 369
$ws.Range("C36").Value = 624
$ws.Range("D36").Value = 1000
$ws.Range("E36").Value = 1200
$ws.Range("F36").Value = 189000
$ws.Range("G36").Value = 111600
$ws.Range("B37").Value = 409
$ws.Range("C37").Value = 358
$ws.Range("D37").Value = 400
$ws.Range("E37").Value = 400
$ws.Range("F37").Value = 136400
$ws.Range("G37").Value = 220800
$ws.Range("B38").Value = 313
$ws.Range("C38").Value = 428
$ws.Range("D38").Value = 700
$ws.Range("E38").Value = 400
$ws.Range("F38").Value = 216300
$ws.Range("G38").Value = 142400
$ws.Range("B39").Value = 219
$ws.Range("C39").Value = 467
$ws.Range("D39").Value = 400
$ws.Range("E39").Value = 800
$ws.Range("F39").Value = 118000
$ws.Range("G39").Value = 222400
$ws.Range("B40").Value = 395
$ws.Range("C40").Value = 471
$ws.Range("D40").Value = 1300
$ws.Range("E40").Value = 400
$ws.Range("F40").Value = 120900
$ws.Range("G40").Value = 186000
$ws.Range("B41").Value = 487
$ws.Range("C41").Value = 520
$ws.Range("D41").Value = 500
$ws.Range("E41").Value = 600
$ws.Range("F41").Value = 152500
$ws.Range("G41").Value = 186000
$ws.Range("B42").Value = 567
$ws.Range("C42").Value = 243
$ws.Range("D42").Value = 400
$ws.Range("E42").Value = 400
$ws.Range("F42").Value = 154400
$ws.Range("G42").Value = 201600
$ws.Range("B43").Value = 479
$ws.Range("C43").Value = 589
$ws.Range("D43").Value = 400
$ws.Range("E43").Value = 1100
$ws.Range("F43").Value = 222800
$ws.Range("G43").Value = 154000
$ws.Range("B44").Value = 386
$ws.Range("C44").Value = 283
$ws.Range("D44").Value = 400
$ws.Range("E44").Value = 400
$ws.Range("F44").Value = 165200
$ws.Range("G44").Value = 213200
$ws.Range("B45").Value = 390
$ws.Range("C45").Value = 392
$ws.Range("D45").Value = 600
$ws.Range("E45").Value = 700
$ws.Range("F45").Value = 211800
$ws.Range("G45").Value = 142800
$ws.Range("B46").Value = 461
$ws.Range("C46").Value = 525
$ws.Range("D46").Value = 500
$ws.Range("E46").Value = 600
$ws.Range("F46").Value = 170000
$ws.Range("G46").Value = 212400
$ws.Range("B47").Value = 444
$ws.Range("C47").Value = 426
$ws.Range("D47").Value = 400
$ws.Range("E47").Value = 500
$ws.Range("F47").Value = 134400
$ws.Range("G47").Value = 176500
$ws.Range("B48").Value = 535
$ws.Range("C48").Value = 444
$ws.Range("D48").Value = 500
$ws.Range("E48").Value = 500
$ws.Range("F48").Value = 160500
$ws.Range("G48").Value = 158000
$ws.Range("B49").Value = 249
$ws.Range("C49").Value = 496
$ws.Range("D49").Value = 300
$ws.Range("E49").Value = 600
$ws.Range("F49").Value = 147900
$ws.Range("G49").Value = 205200
$ws.Range("B50").Value = 423
$ws.Range("C50").Value = 271
$ws.Range("D50").Value = 800
$ws.Range("E50").Value = 400
$ws.Range("F50").Value = 134400
$ws.Range("G50").Value = 194000
$ws.Range("B51").Value = 346
$ws.Range("C51").Value = 399
$ws.Range("D51").Value = 500
$ws.Range("E51").Value = 1000
$ws.Range("F51").Value = 193000
$ws.Range("G51").Value = 138000
